# "new runsettings, testing CT5 examples, working"
#
# Adds two new worksheets to the workbook:
#   - "zero" : inserted as the very first sheet, all exposure/expense
#              figures zeroed out (template/"zero" run settings sheet).
#   - "ea1"  : inserted as the very last sheet, holding a new example
#              (CT5) data set, and left as the active sheet/tab.
#
# Both new sheets are built by copying the existing "ufs1" sheet (so the
# column widths / header styling / page setup match the rest of the
# workbook) and then overwriting the copied data.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("ufs1")

# ---------------------------------------------------------------------
# 1. "zero" sheet - inserted before "ufs1" (i.e. becomes the first tab)
# ---------------------------------------------------------------------
$template.Copy($template)
$zero = $wb.Worksheets.Item(1)
$zero.Name = "zero"

# Drop the inherited numeric-format style on column G so the copied
# cells fall back to the default style (matches a freshly authored
# sheet rather than a clone of "ufs1").
$zero.Range("G2:G6").ClearFormats()

# Zero every data column (SPCODE in column A is left untouched).
$zero.Range("B2:I6").Value = 0

$zero.Range("C9").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. "ea1" sheet - inserted after "ulp" (i.e. becomes the last tab)
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$template.Copy($null, $lastSheet)
$ea1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ea1.Name = "ea1"

$ea1.Range("G2:G6").ClearFormats()

# Row 2 (SPCODE 10)
$ea1.Range("B2").Value = 325
$ea1.Range("C2").Value = 0
$ea1.Range("D2").Value = 1
$ea1.Range("E2").Formula = "=75/12"
$ea1.Range("F2").Value = 0
$ea1.Range("G2").Value = 0.025
$ea1.Range("H2").Value = 0
$ea1.Range("I2").Value = 0

# Row 3 (SPCODE 11)
$ea1.Range("B3:I3").Value = 0

# Row 4 (SPCODE 20)
$ea1.Range("B4:I4").Value = 0

# Row 5 (SPCODE 30)
$ea1.Range("B5:I5").Value = 0

# Row 6 (SPCODE 40)
$ea1.Range("B6").Value = 150
$ea1.Range("C6").Value = 0
$ea1.Range("D6").Value = 0
$ea1.Range("E6").Value = 0
$ea1.Range("F6").Value = 0
$ea1.Range("G6").Value = 0
$ea1.Range("H6").Value = 0
$ea1.Range("I6").Value = 0.03

$ea1.Range("E3").Select() | Out-Null

# "ea1" is the sheet that is on-screen/active when the workbook is saved.
$ea1.Activate()
